# Updated cryptos list values (Price / Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. D-column values that look like plain numbers are
# prefixed with a single quote so Excel keeps them as literal text (matching
# the original inlineStr cells, e.g. "0.410" must not collapse to 0.41).
$updates = [ordered]@{
    "D2" = "26.661.28"
    "E2" = "  +1.07%  "
    "D3" = "1.632.67"
    "E3" = "  +0.45%  "
    "E4" = "  +0.40%  "
    "D5" = "'213.61"
    "E5" = "  +0.47%  "
    "E6" = "  +2.92%  "
    "E7" = "  +0.41%  "
    "D8" = "'0.252"
    "E8" = "  +1.03%  "
    "E9" = "  +0.82%  "
    "D10" = "'19.18"
    "E10" = "  +1.51%  "
    "E11" = "  +3.48%  "
    "E12" = "  +0.54%  "
    "D13" = "1.644.26"
    "E13" = "  +1.14%  "
    "D14" = "'4.10"
    "E14" = "  +1.43%  "
    "E15" = "  +0.23%  "
    "D16" = "26.663.83"
    "E16" = "  +1.04%  "
    "D17" = "'63.57"
    "E17" = "  +1.47%  "
    "D18" = "0.0₃0741"
    "E18" = "  +1.80%  "
    "D19" = "'215.50"
    "E19" = "  +6.22%  "
    "E20" = "  +0.38%  "
    "E21" = "  +0.96%  "
    "E22" = "  +1.67%  "
    "D23" = "'9.35"
    "E23" = "  +0.28%  "
    "D24" = "'1.96"
    "E24" = "  +5.07%  "
    "D25" = "'147.79"
    "E25" = "  +2.17%  "
    "D26" = "'1.01"
    "E26" = "  +0.37%  "
    "E27" = "  +1.38%  "
    "E28" = "  +3.95%  "
    "D29" = "'15.54"
    "E29" = "  +2.19%  "
    "E30" = "  -2.68%  "
    "E31" = "  -0.13%  "
    "E32" = "  +3.21%  "
    "E33" = "  +1.96%  "
    "E34" = "  +0.42%  "
    "D35" = "1.227.26"
    "E36" = "  +0.09%  "
    "D37" = "'0.0172"
    "E37" = "  +5.79%  "
    "E38" = "  -0.12%  "
    "D39" = "'1.01"
    "E39" = "  +0.44%  "
    "E40" = "  +0.94%  "
    "D41" = "'2.29"
    "E41" = "  -1.50%  "
    "D42" = "'0.796"
    "E42" = "  +1.68%  "
    "D43" = "'5.35"
    "E43" = "  -0.74%  "
    "D44" = "1.769.54"
    "E44" = "  +0.30%  "
    "E45" = "  +1.02%  "
    "D46" = "'1.56"
    "E46" = "  +2.60%  "
    "D47" = "'55.26"
    "E47" = "  +2.28%  "
    "D48" = "0.0₆0104"
    "E48" = "  +0.16%  "
    "E49" = "  +1.15%  "
    "D50" = "'7.60"
    "E50" = "  +3.50%  "
    "D51" = "'0.410"
    "E51" = "  +0.15%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
